$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rotate data among rows 2, 3, 4 for columns D, J, K, L, M, N, O, P, Q
# old row2 -> new row4
# old row3 -> new row2
# old row4 -> new row3

$oldD2 = $ws.Cells.Item(2, 4).Value2
$oldJ2 = $ws.Cells.Item(2, 10).Value2
$oldK2 = $ws.Cells.Item(2, 11).Value2
$oldL2 = $ws.Cells.Item(2, 12).Value2
$oldM2 = $ws.Cells.Item(2, 13).Value2
$oldN2 = $ws.Cells.Item(2, 14).Value2
$oldO2 = $ws.Cells.Item(2, 15).Value2
$oldP2 = $ws.Cells.Item(2, 16).Value2
$oldQ2 = $ws.Cells.Item(2, 17).Value2

$oldD3 = $ws.Cells.Item(3, 4).Value2
$oldJ3 = $ws.Cells.Item(3, 10).Value2
$oldK3 = $ws.Cells.Item(3, 11).Value2
$oldL3 = $ws.Cells.Item(3, 12).Value2
$oldM3 = $ws.Cells.Item(3, 13).Value2
$oldN3 = $ws.Cells.Item(3, 14).Value2
$oldO3 = $ws.Cells.Item(3, 15).Value2
$oldP3 = $ws.Cells.Item(3, 16).Value2
$oldQ3 = $ws.Cells.Item(3, 17).Value2

$oldD4 = $ws.Cells.Item(4, 4).Value2
$oldJ4 = $ws.Cells.Item(4, 10).Value2
$oldK4 = $ws.Cells.Item(4, 11).Value2
$oldL4 = $ws.Cells.Item(4, 12).Value2
$oldM4 = $ws.Cells.Item(4, 13).Value2
$oldN4 = $ws.Cells.Item(4, 14).Value2
$oldO4 = $ws.Cells.Item(4, 15).Value2
$oldP4 = $ws.Cells.Item(4, 16).Value2
$oldQ4 = $ws.Cells.Item(4, 17).Value2

# New row 2 = old row 3
$ws.Cells.Item(2, 4).Value2 = $oldD3
$ws.Cells.Item(2, 10).Value2 = $oldJ3
$ws.Cells.Item(2, 11).Value2 = $oldK3
$ws.Cells.Item(2, 12).Value2 = $oldL3
$ws.Cells.Item(2, 13).Value2 = $oldM3
$ws.Cells.Item(2, 14).Value2 = $oldN3
$ws.Cells.Item(2, 15).Value2 = $oldO3
$ws.Cells.Item(2, 16).Value2 = $oldP3
$ws.Cells.Item(2, 17).Value2 = $oldQ3

# New row 3 = old row 4
$ws.Cells.Item(3, 4).Value2 = $oldD4
$ws.Cells.Item(3, 10).Value2 = $oldJ4
$ws.Cells.Item(3, 11).Value2 = $oldK4
$ws.Cells.Item(3, 12).Value2 = $oldL4
$ws.Cells.Item(3, 13).Value2 = $oldM4
$ws.Cells.Item(3, 14).Value2 = $oldN4
$ws.Cells.Item(3, 15).Value2 = $oldO4
$ws.Cells.Item(3, 16).Value2 = $oldP4
$ws.Cells.Item(3, 17).Value2 = $oldQ4

# New row 4 = old row 2
$ws.Cells.Item(4, 4).Value2 = $oldD2
$ws.Cells.Item(4, 10).Value2 = $oldJ2
$ws.Cells.Item(4, 11).Value2 = $oldK2
$ws.Cells.Item(4, 12).Value2 = $oldL2
$ws.Cells.Item(4, 13).Value2 = $oldM2
$ws.Cells.Item(4, 14).Value2 = $oldN2
$ws.Cells.Item(4, 15).Value2 = $oldO2
$ws.Cells.Item(4, 16).Value2 = $oldP2
$ws.Cells.Item(4, 17).Value2 = $oldQ2

$wb.Save()
